# Prior-to-refactor update of the MLR results workbook:
#   - the weather regressions now also control for trip-level covariates
#     (avg_trip_miles, avg_trip_time_min, demand_resid, driver_pay_pct_of_base_fare)
#     in addition to the rain/weather terms, so the model-formula labels on the
#     Summary sheet and the tab names of the per-model sheets need to change.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the per-model worksheet tabs to reflect the new formulas ---
$ws = $wb.Worksheets.Item("fare_per_mile_resid ~ rain_flag")
$ws.Name = "fare_per_mile_resid ~ avg_trip_"

$ws = $wb.Worksheets.Item("margin_per_mile_resid ~ rain_fl")
$ws.Name = "margin_per_mile_resid ~ avg_tri"

$ws = $wb.Worksheets.Item("driverpay_per_mile_resid ~ rain")
$ws.Name = "driverpay_per_mile_resid ~ avg_"

# --- 2. Update the model_label column (A) on the Summary sheet with the new formulas ---
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A2").Value = "avg_base_passenger_fare_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare"
$summary.Range("A3").Value = "fare_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare"
$summary.Range("A4").Value = "margin_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare"

$summary.Range("A5").Value = "avg_base_passenger_fare_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A6").Value = "fare_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A7").Value = "margin_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A8").Value = "fare_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A9").Value = "driverpay_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A10").Value = "margin_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + driver_pay_pct_of_base_fare + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A11").Value = "driverpay_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f"
$summary.Range("A12").Value = "driverpay_per_mile_resid ~ avg_trip_miles + avg_trip_time_min + demand_resid + rain_flag_lag0 + heavy_rain_flag_lag0 + precip_1h_mm_total + wind_chill_f + driver_pay_pct_of_base_fare"
